# Addition of case study examples.
# - Adds a new "Runs" worksheet summarising a curated set of TCR sequences.
# - Adds a new "AS01" TCR example (EBV / BMLF1_280-288) to the "TCRs" sheet.
# - Adds the two corresponding sequence rows (AS01_1 / AS01_2) to "SeqInfo".
# - Highlights the TCR_ID cells for TCRs that have associated "Runs" (yellow fill).

$wb = $excel.ActiveWorkbook

$tcrs = $wb.Worksheets.Item("TCRs")
$seqInfo = $wb.Worksheets.Item("SeqInfo")

# ---------------------------------------------------------------------------
# 1. Sequence IDs for the new AS01 entries (written first).
# ---------------------------------------------------------------------------
$seqInfo.Range("A16").Value = "AS01_1"
$seqInfo.Range("A17").Value = "AS01_2"

# ---------------------------------------------------------------------------
# 2. TCRs sheet: new row 9 ("AS01").
# ---------------------------------------------------------------------------
$tcrs.Range("A9").Value = "AS01"
$tcrs.Range("D9").Value = "GLCTLVAML"
$tcrs.Range("H9").Value = "20-1"
$tcrs.Range("F9").Value = "5"
$tcrs.Range("G9").Value = "31"
$tcrs.Range("C9").Value = "BMLF1_280_288"
$tcrs.Range("J9").Value = "Equal Bias: https://journals.plos.org/plospathogens/article?id=10.1371/journal.ppat.1001198"
$tcrs.Range("B9").Value = "EBV"
$tcrs.Range("E9").Value = "HLA-A*0201"
$tcrs.Range("I9").Value = "1-2"

# Yellow-highlight the TCR_ID column for rows that now have run data
# (JM22 row 2, YVL_TCR row 8, AS01 row 9).
$tcrs.Range("A2").Interior.Color = 65535
$tcrs.Range("A8").Interior.Color = 65535
$tcrs.Range("A9").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 3. SeqInfo sheet: finish new rows 16 & 17 for the AS01 TCR's two sequences.
# ---------------------------------------------------------------------------
$seqInfo.Range("C16").Value = "CAEDFNARLMF"
$seqInfo.Range("D16").Value = "TRAJ31*01"
$seqInfo.Range("E16").Value = "TRBV20-1*01"
$seqInfo.Range("F16").Value = "CSARTGVGNTIYF"
$seqInfo.Range("G16").Value = "TRBJ1-3*01"
$seqInfo.Range("B16").Value = "TRAV5*01"
$seqInfo.Range("H16").Value = 1543

$seqInfo.Range("C17").Value = "CAEDKDARLMF"
$seqInfo.Range("F17").Value = "CSARDRIGNTIYF"
$seqInfo.Range("B17").Value = "TRAV5*01"
$seqInfo.Range("D17").Value = "TRAJ31*01"
$seqInfo.Range("E17").Value = "TRBV20-1*01"
$seqInfo.Range("G17").Value = "TRBJ1-3*01"
$seqInfo.Range("H17").Value = 1546

# ---------------------------------------------------------------------------
# 4. New "Runs" worksheet: curated subset of SeqInfo (JM22, YVL, AS01) with a
#    run/read-count column.
# ---------------------------------------------------------------------------
$runs = $wb.Worksheets.Add()
$runs.Name = "Runs"
$runs.Move([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-fetch a fresh handle (Move() can stale the previous reference).
$runs = $wb.Worksheets.Item("Runs")

$runs.Range("A1").Value = "SeqID"
$runs.Range("B1").Value = "TRAV"
$runs.Range("C1").Value = "CDR3A"
$runs.Range("D1").Value = "TRAJ"
$runs.Range("E1").Value = "TRBV"
$runs.Range("F1").Value = "CDR3B"
$runs.Range("G1").Value = "TRBJ"

$runs.Range("A2").Value = "JM22_1"
$runs.Range("B2").Value = "TRAV27*01"
$runs.Range("C2").Value = "CAGAGSQGNLIF"
$runs.Range("D2").Value = "TRAJ42*01"
$runs.Range("E2").Value = "TRBV19*01"
$runs.Range("F2").Value = "CASSSRSSYEQYF"
$runs.Range("G2").Value = "TRBJ2-7*01"
$runs.Range("H2").Value = 114

$runs.Range("A3").Value = "JM22_2"
$runs.Range("B3").Value = "TRAV27*01"
$runs.Range("C3").Value = "CAGAIGPSNTGKLIF"
$runs.Range("D3").Value = "TRAJ37*01"
$runs.Range("E3").Value = "TRBV19*01"
$runs.Range("F3").Value = "CASSIRSSYEQYF"
$runs.Range("G3").Value = "TRBJ2-7*01"
$runs.Range("H3").Value = 259

$runs.Range("A4").Value = "YVL_1"
$runs.Range("B4").Value = "TRAV13-1*01"
$runs.Range("C4").Value = "CAVKDTDKLIF"
$runs.Range("D4").Value = "TRAJ34*01"
$runs.Range("E4").Value = "TRBV19*01"
$runs.Range("F4").Value = "MSLLGSNQPQHF"
$runs.Range("G4").Value = "TRBJ1-5*01"
$runs.Range("H4").Value = 67

$runs.Range("A5").Value = "YVL_2"
$runs.Range("B5").Value = "TRAV13-1*01"
$runs.Range("C5").Value = "CAVKDTDKLIF"
$runs.Range("D5").Value = "TRAJ34*01"
$runs.Range("E5").Value = "TRBV19*01"
$runs.Range("F5").Value = "CASTGGPGYGAQYF"
$runs.Range("G5").Value = "TRBJ2-5*01"
$runs.Range("H5").Value = 68

$runs.Range("A6").Value = "AS01_1"
$runs.Range("B6").Value = "TRAV5*01"
$runs.Range("C6").Value = "CAEDFNARLMF"
$runs.Range("D6").Value = "TRAJ31*01"
$runs.Range("E6").Value = "TRBV20-1*01"
$runs.Range("F6").Value = "CSARTGVGNTIYF"
$runs.Range("G6").Value = "TRBJ1-3*01"
$runs.Range("H6").Value = 1543

$runs.Range("A7").Value = "AS01_2"
$runs.Range("B7").Value = "TRAV5*01"
$runs.Range("C7").Value = "CAEDKDARLMF"
$runs.Range("D7").Value = "TRAJ31*01"
$runs.Range("E7").Value = "TRBV20-1*01"
$runs.Range("F7").Value = "CSARDRIGNTIYF"
$runs.Range("G7").Value = "TRBJ1-3*01"
$runs.Range("H7").Value = 1546

# ---------------------------------------------------------------------------
# 5. Selections / active sheet, matching the end-state of the author's edit.
# ---------------------------------------------------------------------------
$tcrs.Range("D9").Select()
$seqInfo.Range("A16:H17").Select()
$runs.Range("F10").Select()

# "Runs" (the newly-added, third sheet) ends up the active tab.
$runs.Activate()
